$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.950.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.905.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8006"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.98"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3127"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.31"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07092"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07979"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.905.38"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7386"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.182"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.60"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.962.33"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.96"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.877"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.09"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007785"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.150.97"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9994"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.915"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.48"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.198"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1419"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +9.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.86"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.044"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.360"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.514"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.293"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05561"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.062"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.266"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7301"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.715"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.785"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4407"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.003"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.15"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9993"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8370"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.873"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.572"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.45"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.751"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "976.53"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +7.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.059.49"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.22"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.14%  "
